# Update the table style used by the three data tables in this deck from the
# custom "Table_0" style ({9DFDF103-63D8-4BF3-8B84-76CF84215F48}) to the
# built-in table style {7404B150-1D22-4BE6-82D9-9A7D74A88613}.
$p = $ppt.ActivePresentation

$oldStyleId = "{9DFDF103-63D8-4BF3-8B84-76CF84215F48}"
$newStyleId = "{7404B150-1D22-4BE6-82D9-9A7D74A88613}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
